# chore: adapt column header formatting to respective input file names (#7)
#
# 1. Rename the two sets of column headers (row 1) from the generic
#    "_old" / "_new" suffixes to the concrete format-version suffixes
#    "_FV2210" / "_FV2304".
# 2. Turn the used range (A1:U87) into a real Excel Table ("Table1")
#    with those headers, so the sheet gets the `xl/tables/table1.xml`
#    part + the `<tableParts>` wiring.
# 3. Freeze the header row (split at row 2 / top-left A2, frozen pane)
#    so the header stays visible while scrolling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header row -------------------------------------------------

$headers = @(
    "Segmentname_FV2210",
    "Segmentgruppe_FV2210",
    "Segment_FV2210",
    "Datenelement_FV2210",
    "Segment ID_FV2210",
    "Code_FV2210",
    "Qualifier_FV2210",
    "Beschreibung_FV2210",
    "Bedingungsausdruck_FV2210",
    "Bedingung_FV2210",
    "diff",
    "Segmentname_FV2304",
    "Segmentgruppe_FV2304",
    "Segment_FV2304",
    "Datenelement_FV2304",
    "Segment ID_FV2304",
    "Code_FV2304",
    "Qualifier_FV2304",
    "Beschreibung_FV2304",
    "Bedingungsausdruck_FV2304",
    "Bedingung_FV2304"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value2 = $headers[$i]
}

# --- 2. Convert the range into a Table ------------------------------------

$dataRange = $ws.Range("A1:U87")
$table = $ws.ListObjects.Add(1, $dataRange, [System.Type]::Missing, 1)
$table.Name = "Table1"

# --- 3. Freeze the header row ---------------------------------------------

$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
